$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price/volume/name data (GitHub Actions scheduled update).
# All Coin/Link/Price/Volume cells are stored as plain text in this sheet
# (prices mix thousands-dot formatting, leading zeros, subscript digits, etc,
# so they are never real numbers). Force NumberFormat to Text ("@") before
# writing so Excel does not "helpfully" reinterpret e.g. "0.370" / "20.20" as
# numbers and drop the trailing zero.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.405.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.986.33"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.37"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -8.69%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.24%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.35"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.370"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.76"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0746"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0987"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.95%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.277.84"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.20"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.753"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.04"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.983.52"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.373.32"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.67"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0801"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.25"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "221.20"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.39"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -9.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.65"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.128"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.82"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.33"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.117"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0604"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.25"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.34"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.01"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.452.64"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0919"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0201"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.09%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.10"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -9.30%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.34"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "14.94"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.989"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.68"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.23%  "
